$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Column constants -------------------------------------------------
$colMode      = 1
$colReferent  = 2
$colDirigent  = 3
$colDato      = 4

# --- Simple per-cell name updates (rows keep their identity) ---------
# Row for "Mode 3" -> Dirigent: Mathis -> Mathias O (with slashed O)
$t.Cell(4, $colDirigent).Range.Text  = "Mathias " + [char]0x00D8

# Row for "Mode 4" -> Dirigent: Mathias O -> Kristoffer
$t.Cell(5, $colDirigent).Range.Text  = "Kristoffer"

# Row for "Mode 5" -> Dirigent: Kristoffer -> Mikkel
$t.Cell(6, $colDirigent).Range.Text  = "Mikkel"

# Row for "Mode 6" -> Referent: Mathis -> Mathias O ; Dirigent: Mikkel -> Mathias N
$t.Cell(7, $colReferent).Range.Text  = "Mathias " + [char]0x00D8
$t.Cell(7, $colDirigent).Range.Text  = "Mathias N"

# Row for "Mode 7" -> Referent: Mathias O -> Kristoffer ; Dirigent: Mathias N -> Rasmus
$t.Cell(8, $colReferent).Range.Text  = "Kristoffer"
$t.Cell(8, $colDirigent).Range.Text  = "Rasmus"

# Row for "Mode 8" -> Referent: Kristoffer -> Mikkel ; Dirigent: Rasmus -> Mathias O
$t.Cell(9, $colReferent).Range.Text  = "Mikkel"
$t.Cell(9, $colDirigent).Range.Text  = "Mathias " + [char]0x00D8

# Remove the "_GoBack" bookmark anchored on the "8" run in the Mode cell
# of row 9 (table row index 9 = meeting "8"); it is deleted together with
# its matching bookmarkEnd once we delete the whole "Mode 9" row below,
# which is where Word had placed the bookmarkEnd marker.

# Row that was "Mode 10" -> Dirigent: Mathias O -> Kristoffer (Referent Rasmus unchanged)
$t.Cell(11, $colDirigent).Range.Text = "Kristoffer"

# Row that was "Mode 11" -> Dirigent: Kristoffer -> Mikkel (Referent Mathias N unchanged,
# this absorbs the value that used to belong to the now-deleted "Mode 12" row)
$t.Cell(12, $colDirigent).Range.Text = "Mikkel"

# --- Drop the two rows that corresponded to Mathis being removed -----
# Delete from the bottom up so earlier row indices stay valid.
$t.Rows.Item(13).Delete()   # old "Mode 12" row (Mathis / Mikkel)
$t.Rows.Item(10).Delete()   # old "Mode 9" row  (Mikkel / Mathis)

# --- Renumber the meetings that shifted up after the deletions -------
$t.Cell(10, $colMode).Range.Text = "9"
$t.Cell(11, $colMode).Range.Text = "10"

# --- Replace the trailing empty paragraph with one carrying the
#     "_GoBack" bookmark (this is where Word leaves it after the last
#     edit point in the document) ---------------------------------------
$last = $d.Paragraphs.Last.Range
$last.Collapse(1)  # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $last)
